$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3, column C: "3 heures" becomes "2 heures"
$ws.Range("C3").Value2 = "2 heures"

# New row 4: same date format as rows above (copy formats from A3), new task description, "2 heures"
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value2 = 43228

$ws.Range("B4").Value2 = "Création des uses cases et scénarios "
$ws.Range("C4").Value2 = "2 heures"

$ws.Range("C4").Select()
$excel.ActiveWindow.Zoom = 125
